# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the numeric/percent-looking Price & Volume(1h) cells
# so Excel does not silently coerce them into numbers (matches the inline
# string type used by the rest of the sheet).
$deCells = @("D2","E2","D3","E3","D4","E4","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","E24","D25","E25","D26","E26","D27","D28","E28","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","E46","E47","E48","D49","E49","D50","E50")
foreach ($ref in $deCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Coin / Link column updates (plain text, no coercion risk)
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B21").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C21").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"

# Price / Volume(1h) column updates
$ws.Range("D2").Value = "244.40"
$ws.Range("E2").Value = "-1.18%"
$ws.Range("D3").Value = "27.45"
$ws.Range("E3").Value = "3.82%"
$ws.Range("D4").Value = "5.046"
$ws.Range("E4").Value = "-0.77%"
$ws.Range("E5").Value = "1.08%"
$ws.Range("D6").Value = "6.476"
$ws.Range("E6").Value = "-0.56%"
$ws.Range("D7").Value = "0.8222"
$ws.Range("E7").Value = "1.10%"
$ws.Range("D8").Value = "0.8394"
$ws.Range("E8").Value = "-1.32%"
$ws.Range("D9").Value = "0.0006012"
$ws.Range("E9").Value = "0.73%"
$ws.Range("D10").Value = "0.1326"
$ws.Range("E10").Value = "-1.30%"
$ws.Range("D11").Value = "0.06934"
$ws.Range("E11").Value = "-0.47%"
$ws.Range("D12").Value = "0.02861"
$ws.Range("E12").Value = "1.24%"
$ws.Range("D13").Value = "0.09388"
$ws.Range("E13").Value = "-0.21%"
$ws.Range("D14").Value = "0.001511"
$ws.Range("E14").Value = "-1.05%"
$ws.Range("D15").Value = "0.04134"
$ws.Range("E15").Value = "-11.29%"
$ws.Range("D16").Value = "0.006119"
$ws.Range("E16").Value = "-1.95%"
$ws.Range("D17").Value = "3.509"
$ws.Range("E17").Value = "-2.13%"
$ws.Range("D18").Value = "3.001"
$ws.Range("E18").Value = "-1.60%"
$ws.Range("D19").Value = "2.307"
$ws.Range("E19").Value = "8.89%"
$ws.Range("D20").Value = "0.3113"
$ws.Range("E20").Value = "-2.15%"
$ws.Range("D21").Value = "0.03160"
$ws.Range("E21").Value = "-0.24%"
$ws.Range("D22").Value = "0.1254"
$ws.Range("E22").Value = "-4.98%"
$ws.Range("D23").Value = "3.582"
$ws.Range("E23").Value = "-4.25%"
$ws.Range("E24").Value = "1.81%"
$ws.Range("D25").Value = "0.001221"
$ws.Range("E25").Value = "-2.19%"
$ws.Range("D26").Value = "0.003869"
$ws.Range("E26").Value = "-16.18%"
$ws.Range("D27").Value = "0.00009804"
$ws.Range("D28").Value = "0.0001439"
$ws.Range("E28").Value = "-25.72%"
$ws.Range("D40").Value = "0.03698"
$ws.Range("E40").Value = "0.57%"
$ws.Range("D41").Value = "0.006153"
$ws.Range("E41").Value = "80.35%"
$ws.Range("D42").Value = "0.1053"
$ws.Range("E42").Value = "-22.09%"
$ws.Range("D43").Value = "0.002294"
$ws.Range("E43").Value = "-13.77%"
$ws.Range("D44").Value = "0.009561"
$ws.Range("E44").Value = "11.53%"
$ws.Range("D45").Value = "0.00005194"
$ws.Range("E45").Value = "-1.80%"
$ws.Range("E46").Value = "0.08%"
$ws.Range("E47").Value = "-15.36%"
$ws.Range("E48").Value = "0.08%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "0.08%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "0.08%"

# Strip the temporary text format so the cells fall back to the default,
# unstyled state (matching the rest of the sheet) while keeping their text value.
foreach ($ref in $deCells) {
    $ws.Range($ref).ClearFormats()
}
